# Apply the deck-wide date placeholder refresh (10/30/2017 -> 10/31/2017)
# across the slide master and every custom (slide) layout, then update the
# title text on the first slide.

$p = $ppt.ActivePresentation

# --- 1. Slide master's Date Placeholder ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "10/30/2017") {
            $sh.TextFrame.TextRange.Text = "10/31/2017"
        }
    }
}

# --- 2. Every custom layout's Date Placeholder ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "10/30/2017") {
                $sh.TextFrame.TextRange.Text = "10/31/2017"
            }
        }
    }
}

# --- 3. Title text on slide 1 ---
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Reproducible Research with R, The Tidyverse, Notebooks, and Spark"
